$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.419.40"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.672.33"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'603.09"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'178.27"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.671.35"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("D14").Value = "3.159.10"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "72.293.54"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "'26.29"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "2.671.57"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "'11.92"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("D20").Value = "'7.98"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'371.28"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +8.96%  "
$ws.Range("D24").Value = "'72.38"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "'9.87"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "2.809.34"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.08"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'517.91"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'164.67"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'19.44"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "'19.13"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -7.43%  "
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'5.03"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'39.22"
$ws.Range("D47").Value = "'152.90"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value = "'0.0767"
$ws.Range("E51").Value = "  +1.85%  "
